$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("January")
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C33").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("February")
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C33").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("March")
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").Value = 0.1458333333333333
$ws.Range("C22").Value = 0.04166666666666666
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.04166666666666666
$ws.Range("C25").Value = 0.1666666666666667
$ws.Range("C26").Value = 0.04166666666666666
$ws.Range("C27").Value = 0.1666666666666667
$ws.Range("C28").Value = 0.1666666666666667
$ws.Range("C29").Value = 0.04166666666666666
$ws.Range("C30").Value = 0.04166666666666666
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.1666666666666667
$ws.Range("C33").Value = 0.05208333333333334
$ws.Range("C34").Value = 0.04166666666666666
$ws.Range("C35").Value = 0.0625
$ws.Range("C36").Value = 0.0625
$ws.Range("C37").Value = 0.15625

$ws = $wb.Worksheets.Item("April")
$ws.Range("C7").Value = 0.09375
$ws.Range("C8").Value = 0.04166666666666666
$ws.Range("C9").Value = 0.04166666666666666
$ws.Range("C10").Value = 0.04166666666666666
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.04166666666666666
$ws.Range("C13").Value = 0.04166666666666666
$ws.Range("C14").Value = 0.1666666666666667
$ws.Range("C15").Value = 0.1666666666666667
$ws.Range("C16").Value = 0.09375
$ws.Range("C17").Value = 0.04166666666666666
$ws.Range("C18").Value = 0.1458333333333333
$ws.Range("C19").Value = 0.04166666666666666
$ws.Range("C20").Value = 0.04166666666666666
$ws.Range("C21").Value = 0.05208333333333334
$ws.Range("C22").Value = 0.1458333333333333
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.1666666666666667
$ws.Range("C25").Value = 0.05208333333333334
$ws.Range("C26").Value = 0.04166666666666666
$ws.Range("C27").Value = 0.04166666666666666
$ws.Range("C28").Value = 0.1666666666666667
$ws.Range("C29").Value = 0.125
$ws.Range("C30").Value = 0.04166666666666666
$ws.Range("C31").Value = 0.1666666666666667
$ws.Range("C32").Value = 0.05208333333333334
$ws.Range("C33").Value = 0.1666666666666667
$ws.Range("C34").Value = 0.04166666666666666
$ws.Range("C35").Value = 0.0625
$ws.Range("C36").Value = 0.05208333333333334
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("May")
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").Value = 0.1666666666666667
$ws.Range("C9").Value = 0.08333333333333333
$ws.Range("C10").Value = 0.1145833333333333
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.04166666666666666
$ws.Range("C13").Value = 0.08333333333333333
$ws.Range("C14").Value = 0.04166666666666666
$ws.Range("C15").Value = 0.04166666666666666
$ws.Range("C16").Value = 0.04166666666666666
$ws.Range("C17").Value = 0.08333333333333333
$ws.Range("C18").Value = 0.04166666666666666
$ws.Range("C19").Value = 0.1666666666666667
$ws.Range("C20").Value = 0.04166666666666666
$ws.Range("C21").Value = 0.09375
$ws.Range("C22").Value = 0.04166666666666666
$ws.Range("C23").Value = 0.09375
$ws.Range("C24").Value = 0.04166666666666666
$ws.Range("C25").Value = 0.1666666666666667
$ws.Range("C26").Value = 0.04166666666666666
$ws.Range("C27").Value = 0.04166666666666666
$ws.Range("C28").Value = 0.07291666666666667
$ws.Range("C29").Value = 0.07291666666666667
$ws.Range("C30").Value = 0.08333333333333333
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.04166666666666666
$ws.Range("C33").Value = 0.1666666666666667
$ws.Range("C34").Value = 0.15625
$ws.Range("C35").Value = 0.1666666666666667
$ws.Range("C36").Value = 0.04166666666666666
$ws.Range("C37").Value = 0.1666666666666667

$ws = $wb.Worksheets.Item("June")
$ws.Range("C7").Value = 0.125
$ws.Range("C8").Value = 0.08333333333333333
$ws.Range("C9").Value = 0.1145833333333333
$ws.Range("C10").Value = 0.1145833333333333
$ws.Range("C11").Value = 0.1666666666666667
$ws.Range("C12").Value = 0.1145833333333333
$ws.Range("C13").Value = 0.08333333333333333
$ws.Range("C14").Value = 0.08333333333333333
$ws.Range("C15").Value = 0.04166666666666666
$ws.Range("C16").Value = 0.1666666666666667
$ws.Range("C17").Value = 0.04166666666666666
$ws.Range("C18").Value = 0.04166666666666666
$ws.Range("C19").Value = 0.1145833333333333
$ws.Range("C20").Value = 0.09375
$ws.Range("C21").Value = 0.04166666666666666
$ws.Range("C22").Value = 0.04166666666666666
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.1666666666666667
$ws.Range("C25").Value = 0.1666666666666667
$ws.Range("C26").Value = 0.04166666666666666
$ws.Range("C27").Value = 0.04166666666666666
$ws.Range("C28").Value = 0.04166666666666666
$ws.Range("C29").Value = 0.05208333333333334
$ws.Range("C30").Value = 0.04166666666666666
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.05208333333333334
$ws.Range("C33").Value = 0.04166666666666666
$ws.Range("C34").Value = 0.08333333333333333
$ws.Range("C35").Value = 0.04166666666666666
$ws.Range("C36").Value = 0.04166666666666666
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("July")
$ws.Range("C7").Value = 0.04166666666666666
$ws.Range("C8").Value = 0.0625
$ws.Range("C9").Value = 0.05208333333333334
$ws.Range("C10").Value = 0.04166666666666666
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.04166666666666666
$ws.Range("C13").Value = 0.1666666666666667
$ws.Range("C14").Value = 0.05208333333333334
$ws.Range("C15").Value = 0.09375
$ws.Range("C16").Value = 0.1458333333333333
$ws.Range("C17").Value = 0.1458333333333333
$ws.Range("C18").Value = 0.1666666666666667
$ws.Range("C19").Value = 0.04166666666666666
$ws.Range("C20").Value = 0.05208333333333334
$ws.Range("C21").Value = 0.1666666666666667
$ws.Range("C22").Value = 0.04166666666666666
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.04166666666666666
$ws.Range("C25").Value = 0.05208333333333334
$ws.Range("C26").Value = 0.1666666666666667
$ws.Range("C27").Value = 0.04166666666666666
$ws.Range("C28").Value = 0.1145833333333333
$ws.Range("C29").Value = 0.04166666666666666
$ws.Range("C30").Value = 0.15625
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.04166666666666666
$ws.Range("C33").Value = 0.1666666666666667
$ws.Range("C34").Value = 0.1666666666666667
$ws.Range("C35").Value = 0.04166666666666666
$ws.Range("C36").Value = 0.04166666666666666
$ws.Range("C37").Value = 0.04166666666666666

$ws = $wb.Worksheets.Item("August")
$ws.Range("C7").Value = 0.125
$ws.Range("C8").Value = 0.1666666666666667
$ws.Range("C9").Value = 0.04166666666666666
$ws.Range("C10").Value = 0.1458333333333333
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.1666666666666667
$ws.Range("C13").Value = 0.1666666666666667
$ws.Range("C14").Value = 0.04166666666666666
$ws.Range("C15").Value = 0.04166666666666666
$ws.Range("C16").Value = 0.08333333333333333
$ws.Range("C17").Value = 0.1666666666666667
$ws.Range("C18").Value = 0.04166666666666666
$ws.Range("C19").Value = 0.1041666666666667
$ws.Range("C20").Value = 0.05208333333333334
$ws.Range("C21").Value = 0.04166666666666666
$ws.Range("C22").Value = 0.125
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.04166666666666666
$ws.Range("C25").Value = 0.04166666666666666
$ws.Range("C26").Value = 0.08333333333333333
$ws.Range("C27").Value = 0.1666666666666667
$ws.Range("C28").Value = 0.04166666666666666
$ws.Range("C29").Value = 0.04166666666666666
$ws.Range("C30").Value = 0.07291666666666667
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.1666666666666667
$ws.Range("C33").Value = 0.1145833333333333
$ws.Range("C34").Value = 0.04166666666666666
$ws.Range("C35").Value = 0.1458333333333333
$ws.Range("C36").Value = 0.08333333333333333
$ws.Range("C37").Value = 0.1041666666666667

$ws = $wb.Worksheets.Item("September")
$ws.Range("C7").Value = 0.04166666666666666
$ws.Range("C8").Value = 0.04166666666666666
$ws.Range("C9").Value = 0.1666666666666667
$ws.Range("C10").Value = 0.04166666666666666
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.125
$ws.Range("C13").Value = 0.04166666666666666
$ws.Range("C14").Value = 0.04166666666666666
$ws.Range("C15").Value = 0.09375
$ws.Range("C16").Value = 0.04166666666666666
$ws.Range("C17").Value = 0.04166666666666666
$ws.Range("C18").Value = 0.08333333333333333
$ws.Range("C19").Value = 0.04166666666666666
$ws.Range("C20").Value = 0.04166666666666666
$ws.Range("C21").Value = 0.1666666666666667
$ws.Range("C22").Value = 0.05208333333333334
$ws.Range("C23").Value = 0.04166666666666666
$ws.Range("C24").Value = 0.1458333333333333
$ws.Range("C25").Value = 0.04166666666666666
$ws.Range("C26").Value = 0.07291666666666667
$ws.Range("C27").Value = 0.09375
$ws.Range("C28").Value = 0.04166666666666666
$ws.Range("C29").Value = 0.1666666666666667
$ws.Range("C30").Value = 0.0625
$ws.Range("C31").Value = 0.07291666666666667
$ws.Range("C32").Value = 0.1666666666666667
$ws.Range("C33").Value = 0.1666666666666667
$ws.Range("C34").Value = 0.1145833333333333
$ws.Range("C35").Value = 0.1666666666666667
$ws.Range("C36").Value = 0.1145833333333333
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("October")
$ws.Range("C7").Value = 0.08333333333333333
$ws.Range("C8").Value = 0.09375
$ws.Range("C9").Value = 0.05208333333333334
$ws.Range("C10").Value = 0.04166666666666666
$ws.Range("C11").Value = 0.04166666666666666
$ws.Range("C12").Value = 0.04166666666666666
$ws.Range("C13").Value = 0.1666666666666667
$ws.Range("C14").Value = 0.08333333333333333
$ws.Range("C15").Value = 0.125
$ws.Range("C16").Value = 0.04166666666666666
$ws.Range("C17").Value = 0.07291666666666667
$ws.Range("C18").Value = 0.04166666666666666
$ws.Range("C19").Value = 0.07291666666666667
$ws.Range("C20").Value = 0.0625
$ws.Range("C21").Value = 0.1666666666666667
$ws.Range("C22").Value = 0.05208333333333334
$ws.Range("C23").Value = 0.05208333333333334
$ws.Range("C24").Value = 0.1666666666666667
$ws.Range("C25").Value = 0.1666666666666667
$ws.Range("C26").Value = 0.09375
$ws.Range("C27").Value = 0.0625
$ws.Range("C28").Value = 0.07291666666666667
$ws.Range("C29").Value = 0.04166666666666666
$ws.Range("C30").Value = 0.09375
$ws.Range("C31").Value = 0.04166666666666666
$ws.Range("C32").Value = 0.04166666666666666
$ws.Range("C33").Value = 0.04166666666666666
$ws.Range("C34").Value = 0.04166666666666666
$ws.Range("C35").Value = 0.04166666666666666
$ws.Range("C36").Value = 0.1041666666666667
$ws.Range("C37").Value = 0.1666666666666667

$ws = $wb.Worksheets.Item("November")
$ws.Range("C7").Value = 0.1666666666666667
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C33").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()

$ws = $wb.Worksheets.Item("December")
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("C11").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C20").ClearContents()
$ws.Range("C21").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("C25").ClearContents()
$ws.Range("C26").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").ClearContents()
$ws.Range("C31").ClearContents()
$ws.Range("C32").ClearContents()
$ws.Range("C33").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("C35").ClearContents()
$ws.Range("C36").ClearContents()
$ws.Range("C37").ClearContents()
